# fix dynamic changing of control no
#
# Replaces the sample ICT TA report rows (rows 15 & 16) with the May 2020
# data, blanks out what used to be the third data row (row 17, which is now
# a trailing blank row like row 18), and updates the period title and the
# remembered selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper: write $text into $cellRef as a literal text value (never letting
# Excel's input parser reinterpret an ISO-ish "2020-05-08" string as a real
# date, which would change both the stored type and the cell's style).
# We build the text via a throw-away formula cell (so it is a real string,
# never auto-converted) and PasteSpecial only the *value* over the target,
# leaving the target's existing number format / style index untouched.
function Set-LiteralText($cellRef, $text) {
    $scratch = $ws.Range("Z1")
    $scratch.Formula = "=" + '"' + $text + '"'
    $scratch.Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $scratch.ClearContents()
}

# --- Period title -----------------------------------------------------
$ws.Range("E11").Value = "Month of May 2020"

# --- Row 15 -------------------------------------------------------------
$ws.Range("B15").Value = "2020-243"
Set-LiteralText "C15" "2020-05-08"
$ws.Range("D15").Value = "12:35 PM"
$ws.Range("E15").Value = "NOEL R BARTOLABAC"
$ws.Range("G15").Value = "ORD"
$ws.Range("H15").Value = "ddw"
$ws.Range("I15").Value = "DESKTOP/LAPTOP"
$ws.Range("J15").Value = "Charles Adrian T. Odi"
Set-LiteralText "L15" "2020-05-08"
$ws.Range("M15").Value = "1:09 PM"
$ws.Range("N15").ClearContents()
$ws.Range("O15").Value = "4:28 AM"
$ws.Range("P15").Value = "12 hours and 9 minutes"
$ws.Range("Q15").Value = 0

# --- Row 16 -------------------------------------------------------------
$ws.Range("B16").Value = "2020-244"
Set-LiteralText "C16" "2020-05-13"
$ws.Range("D16").Value = "9:49 AM"
$ws.Range("E16").Value = "NOEL R BARTOLABAC"
$ws.Range("G16").Value = "ORD"
$ws.Range("H16").Value = "test1"
$ws.Range("I16").Value = "SOFTWARE/SYSTEM"
$ws.Range("J16").ClearContents()
$ws.Range("L16").ClearContents()
$ws.Range("M16").Value = "4:28 AM"
$ws.Range("N16").ClearContents()
$ws.Range("O16").Value = "4:28 AM"
$ws.Range("P16").Value = "0 hours and 0 minutes"
$ws.Range("Q16").Value = 0

# --- Row 17: now a blank trailing row, like row 18 ----------------------
$ws.Range("E17:F17").UnMerge()
$ws.Range("A17:Q17").ClearContents()
$ws.Range("A18:Q18").Copy() | Out-Null
$ws.Range("A17:Q17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Rows.Item(17).AutoFit()

# --- Selection ------------------------------------------------------------
$ws.Range("Q16").Select() | Out-Null
